$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("F3").Value = 10.3833
$ws.Range("G3").Value = 46.98

$ws.Range("D5").Value = 11.63715
$ws.Range("E5").Value = 51.39

$ws.Range("F6").Value = 9.8320000000000007
$ws.Range("G6").Value = 45.44

$ws.Range("D7").Value = 12.884499999999999
$ws.Range("E7").Value = 57.52
$ws.Range("F7").Value = 8.6527999999999992
$ws.Range("G7").Value = 38.26

$ws.Range("D11").Value = 13.151999999999999
$ws.Range("E11").Value = 58.1
$ws.Range("F11").Value = 7.9095500000000003
$ws.Range("G11").Value = 32.840000000000003

$ws.Range("D12").Value = 14.7683
$ws.Range("E12").Value = 65.58
$ws.Range("F12").Value = 7.7417999999999996
$ws.Range("G12").Value = 32.89

$ws.Range("D6").Select()
